$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks numeric need to be forced to Text format
# first, otherwise Excel auto-converts them to a numeric cell (losing the
# exact textual representation, e.g. "1.00" -> 1, "0.0000220" -> 2.2E-05).
function Set-TextValue {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
}

$ws.Range("D2").Value = "91.198.38"
$ws.Range("E2").Value = "  +4.07%  "

$ws.Range("D3").Value = "3.084.62"
$ws.Range("E3").Value = "  -0.38%  "

Set-TextValue $ws.Range("D4") "0.998"
$ws.Range("E4").Value = "  -0.27%  "

Set-TextValue $ws.Range("D5") "215.72"
$ws.Range("E5").Value = "  +1.14%  "

Set-TextValue $ws.Range("D6") "618.04"
$ws.Range("E6").Value = "  -2.48%  "

Set-TextValue $ws.Range("D7") "0.374"
$ws.Range("E7").Value = "  -3.08%  "

Set-TextValue $ws.Range("D8") "0.878"
$ws.Range("E8").Value = "  +11.00%  "

Set-TextValue $ws.Range("D9") "1.00"
$ws.Range("E9").Value = "  -0.13%  "

$ws.Range("D10").Value = "3.080.22"
$ws.Range("E10").Value = "  -0.41%  "

Set-TextValue $ws.Range("D11") "0.669"
$ws.Range("E11").Value = "  +19.50%  "

$ws.Range("E12").Value = "  +5.94%  "

$ws.Range("E13").Value = "  +0.22%  "

$ws.Range("D14").Value = "90.824.96"
$ws.Range("E14").Value = "  +3.62%  "

Set-TextValue $ws.Range("D15") "5.36"
$ws.Range("E15").Value = "  +0.01%  "

Set-TextValue $ws.Range("D16") "32.83"
$ws.Range("E16").Value = "  +3.13%  "

$ws.Range("D17").Value = "3.648.15"
$ws.Range("E17").Value = "  -0.71%  "

$ws.Range("D18").Value = "3.078.49"
$ws.Range("E18").Value = "  -1.19%  "

Set-TextValue $ws.Range("D19") "3.49"
$ws.Range("E19").Value = "  +4.30%  "

Set-TextValue $ws.Range("D20") "0.0000220"
$ws.Range("E20").Value = "  +1.27%  "

Set-TextValue $ws.Range("D21") "13.73"
$ws.Range("E21").Value = "  +4.51%  "

Set-TextValue $ws.Range("D22") "433.32"
$ws.Range("E22").Value = "  +3.14%  "

Set-TextValue $ws.Range("D23") "8.45"
$ws.Range("E23").Value = "  +1.05%  "

Set-TextValue $ws.Range("D24") "5.09"
$ws.Range("E24").Value = "  +4.93%  "

Set-TextValue $ws.Range("D25") "5.53"
$ws.Range("E25").Value = "  +2.39%  "

Set-TextValue $ws.Range("D26") "11.84"
$ws.Range("E26").Value = "  +4.26%  "

Set-TextValue $ws.Range("D27") "83.68"
$ws.Range("E27").Value = "  +2.35%  "

$ws.Range("D28").Value = "3.232.36"
$ws.Range("E28").Value = "  -1.56%  "

$ws.Range("E29").Value = "  +0.03%  "

Set-TextValue $ws.Range("D30") "1.08"
$ws.Range("E30").Value = "  +7.68%  "

Set-TextValue $ws.Range("D31") "0.167"
$ws.Range("E31").Value = "  +8.57%  "

Set-TextValue $ws.Range("D32") "8.61"
$ws.Range("E32").Value = "  +6.29%  "

$ws.Range("E33").Value = "  -4.60%  "

Set-TextValue $ws.Range("D34") "514.55"
$ws.Range("E34").Value = "  +3.05%  "

Set-TextValue $ws.Range("D35") "6.85"
$ws.Range("E35").Value = "  +0.32%  "

$ws.Range("E36").Value = "  +0.79%  "

$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D37") "1.26"
$ws.Range("E37").Value = "  +0.09%  "

$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D38") "22.98"
$ws.Range("E38").Value = "  +4.32%  "

$ws.Range("E39").Value = "  -7.14%  "

$ws.Range("E40").Value = "  +0.57%  "

Set-TextValue $ws.Range("D41") "1.00"
$ws.Range("E41").Value = "  -0.35%  "

$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D43") "0.139"
$ws.Range("E43").Value = "  +4.20%  "

$ws.Range("B44").Value = "PolygonEcosystemToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue $ws.Range("D44") "0.367"
$ws.Range("E44").Value = "  +1.51%  "

$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D45") "1.86"
$ws.Range("E45").Value = "  +2.32%  "

Set-TextValue $ws.Range("D46") "0.0718"
$ws.Range("E46").Value = "  +11.08%  "

Set-TextValue $ws.Range("D47") "143.14"
$ws.Range("E47").Value = "  -1.76%  "

Set-TextValue $ws.Range("D48") "0.000270"
$ws.Range("E48").Value = "  +15.00%  "

Set-TextValue $ws.Range("D49") "43.66"
$ws.Range("E49").Value = "  +0.19%  "

Set-TextValue $ws.Range("D50") "4.23"
$ws.Range("E50").Value = "  +8.14%  "

Set-TextValue $ws.Range("D51") "165.74"
$ws.Range("E51").Value = "  +2.91%  "
